$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "Brazilian Serie A"
$ws.Cells.Item(2, 2).NumberFormat = "@"
$ws.Cells.Item(2, 2).Value = "2025-11-25"
$ws.Cells.Item(2, 2).Style = "Normal"
$ws.Cells.Item(2, 3).Value = "21:30:00"
$ws.Cells.Item(2, 4).Value = "Gremio"
$ws.Cells.Item(2, 5).Value = "SE Palmeiras"
$ws.Cells.Item(2, 6).Value = 2.48
$ws.Cells.Item(2, 7).Value = 2.5
$ws.Cells.Item(2, 8).Value = 3.35
$ws.Cells.Item(2, 9).Value = 3.4
$ws.Cells.Item(2, 10).Value = 3.25
$ws.Cells.Item(2, 11).Value = 3.3
$ws.Cells.Item(2, 12).Value = 1.5
$ws.Cells.Item(2, 13).Value = 1.1
$ws.Cells.Item(2, 14).Value = 3.3
$ws.Cells.Item(2, 15).Value = 1.43
$ws.Cells.Item(2, 16).Value = 1.74
$ws.Cells.Item(2, 17).Value = 2.3
$ws.Cells.Item(2, 18).Value = 1.28
$ws.Cells.Item(2, 19).Value = 4.4
$ws.Cells.Item(2, 20).Value = 1.91
$ws.Cells.Item(2, 21).Value = 2.04
$ws.Cells.Item(2, 22).Value = 1.41
$ws.Cells.Item(2, 23).Value = 1.67
$ws.Cells.Item(2, 24).Value = 10.5
$ws.Cells.Item(2, 25).Value = 11
$ws.Cells.Item(2, 26).Value = 22
$ws.Cells.Item(2, 27).Value = 65
$ws.Cells.Item(2, 28).Value = 9
$ws.Cells.Item(2, 29).Value = 7.2
$ws.Cells.Item(2, 30).Value = 14.5
$ws.Cells.Item(2, 31).Value = 36
$ws.Cells.Item(2, 32).Value = 14.5
$ws.Cells.Item(2, 33).Value = 11
$ws.Cells.Item(2, 34).Value = 18.5
$ws.Cells.Item(2, 35).Value = 55
$ws.Cells.Item(2, 36).Value = 38
$ws.Cells.Item(2, 37).Value = 28
$ws.Cells.Item(2, 38).Value = 46
$ws.Cells.Item(2, 39).Value = 130
$ws.Cells.Item(2, 40).Value = 32
$ws.Cells.Item(2, 41).Value = 55

# Row 3
$ws.Cells.Item(3, 1).Value = "Brazilian Serie A"
$ws.Cells.Item(3, 2).NumberFormat = "@"
$ws.Cells.Item(3, 2).Value = "2025-11-25"
$ws.Cells.Item(3, 2).Style = "Normal"
$ws.Cells.Item(3, 3).Value = "21:30:00"
$ws.Cells.Item(3, 4).Value = "Atletico MG"
$ws.Cells.Item(3, 5).Value = "Flamengo"
$ws.Cells.Item(3, 6).Value = 4.4
$ws.Cells.Item(3, 7).Value = 4.5
$ws.Cells.Item(3, 8).Value = 2.08
$ws.Cells.Item(3, 9).Value = 2.1
$ws.Cells.Item(3, 10).Value = 3.3
$ws.Cells.Item(3, 11).Value = 3.35
$ws.Cells.Item(3, 12).Value = 1.55
$ws.Cells.Item(3, 13).Value = 1.11
$ws.Cells.Item(3, 14).Value = 2.92
$ws.Cells.Item(3, 15).Value = 1.5
$ws.Cells.Item(3, 16).Value = 1.62
$ws.Cells.Item(3, 17).Value = 2.56
$ws.Cells.Item(3, 18).Value = 1.23
$ws.Cells.Item(3, 19).Value = 5.2
$ws.Cells.Item(3, 20).Value = 2.12
$ws.Cells.Item(3, 21).Value = 1.85
$ws.Cells.Item(3, 22).Value = 1.9
$ws.Cells.Item(3, 23).Value = 1.28
$ws.Cells.Item(3, 24).Value = 9.4
$ws.Cells.Item(3, 25).Value = 7
$ws.Cells.Item(3, 26).Value = 11
$ws.Cells.Item(3, 27).Value = 28
$ws.Cells.Item(3, 28).Value = 12
$ws.Cells.Item(3, 29).Value = 7
$ws.Cells.Item(3, 30).Value = 11
$ws.Cells.Item(3, 31).Value = 26
$ws.Cells.Item(3, 32).Value = 28
$ws.Cells.Item(3, 33).Value = 17.5
$ws.Cells.Item(3, 34).Value = 23
$ws.Cells.Item(3, 35).Value = 55
$ws.Cells.Item(3, 36).Value = 100
$ws.Cells.Item(3, 37).Value = 70
$ws.Cells.Item(3, 38).Value = 85
$ws.Cells.Item(3, 39).Value = 170
$ws.Cells.Item(3, 40).Value = 100
$ws.Cells.Item(3, 41).Value = 24

# Row 4
$ws.Cells.Item(4, 1).Value = "Colombian Primera A"
$ws.Cells.Item(4, 2).NumberFormat = "@"
$ws.Cells.Item(4, 2).Value = "2025-11-25"
$ws.Cells.Item(4, 2).Style = "Normal"
$ws.Cells.Item(4, 3).Value = "21:30:00"
$ws.Cells.Item(4, 4).Value = "Santa Fe"
$ws.Cells.Item(4, 5).Value = "Tolima"
$ws.Cells.Item(4, 6).Value = 2.48
$ws.Cells.Item(4, 7).Value = 2.52
$ws.Cells.Item(4, 8).Value = 3.55
$ws.Cells.Item(4, 9).Value = 3.65
$ws.Cells.Item(4, 10).Value = 3.05
$ws.Cells.Item(4, 11).Value = 3.2
$ws.Cells.Item(4, 12).Value = 1.54
$ws.Cells.Item(4, 13).Value = 1.12
$ws.Cells.Item(4, 14).Value = 2.94
$ws.Cells.Item(4, 15).Value = 1.48
$ws.Cells.Item(4, 16).Value = 1.62
$ws.Cells.Item(4, 17).Value = 2.48
$ws.Cells.Item(4, 18).Value = 1.24
$ws.Cells.Item(4, 19).Value = 4.9
$ws.Cells.Item(4, 20).Value = 2
$ws.Cells.Item(4, 21).Value = 1.92
$ws.Cells.Item(4, 22).Value = 1.37
$ws.Cells.Item(4, 23).Value = 1.65
$ws.Cells.Item(4, 24).Value = 9.4
$ws.Cells.Item(4, 25).Value = 11
$ws.Cells.Item(4, 26).Value = 24
$ws.Cells.Item(4, 27).Value = 75
$ws.Cells.Item(4, 28).Value = 8.6
$ws.Cells.Item(4, 29).Value = 7
$ws.Cells.Item(4, 30).Value = 16
$ws.Cells.Item(4, 31).Value = 55
$ws.Cells.Item(4, 32).Value = 14
$ws.Cells.Item(4, 33).Value = 11.5
$ws.Cells.Item(4, 34).Value = 19.5
$ws.Cells.Item(4, 35).Value = 100
$ws.Cells.Item(4, 36).Value = 36
$ws.Cells.Item(4, 37).Value = 32
$ws.Cells.Item(4, 38).Value = 75
$ws.Cells.Item(4, 39).Value = 550
$ws.Cells.Item(4, 40).Value = 46
$ws.Cells.Item(4, 41).Value = 980

# Remove rows 5 and 6 (no longer present after the edit)
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(5).Delete()
